$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in A2 from "ZIP Code" to "ZIP_Code"
$ws.Range("A2").Value = "ZIP_Code"

# Move the active selection to A2 (as left by the author after editing it)
$ws.Range("A2").Select()
